{"js": "// Locate the paragraph that begins the Perseus observation instructions\n// (currently split across many runs) and replace its whole content with a\n// single consolidated run of new text, per the commit:\n// \"Print the first paragraph with the name of the constellation.\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst newText =\n  \"Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi \" +\n  \"himmeimpi\u00e4 n\u00e4kyviss\u00e4 olevia t\u00e4hti\u00e4 keinona mitata valonsaastetta tietyss\u00e4 \" +\n  \"paikassa. Paikallistamalla ja tarkkailemalla Perseuksen t\u00e4hdist\u00f6 miten \" +\n  \"valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. \" +\n  \"Antamasi tiedot p\u00e4ivittyv\u00e4t heti verkossa olevaan tietokantaan, ja n\u00e4in \" +\n  \"saadaan k\u00e4sitys siit\u00e4 mink\u00e4 verran taivaan t\u00e4hdist\u00e4 on miss\u00e4kin n\u00e4ht\u00e4viss\u00e4.\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Osallistut maailmanlaajuiseen\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  target.clear();\n  target.insertText(newText, Word.InsertLocation.start);\n  await context.sync();\n}\n", "ps1": "# Locate the paragraph that begins the Perseus observation instructions\n# (currently split across many differently-formatted runs) and replace its\n# whole content with a single consolidated run of new text, per the commit:\n# \"Print the first paragraph with the name of the constellation.\"\n$d = $word.ActiveDocument\n\n$newText = \"Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi himmeimpi\u00e4 n\u00e4kyviss\u00e4 olevia t\u00e4hti\u00e4 keinona mitata valonsaastetta tietyss\u00e4 paikassa. Paikallistamalla ja tarkkailemalla Perseuksen t\u00e4hdist\u00f6 miten valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. Antamasi tiedot p\u00e4ivittyv\u00e4t heti verkossa olevaan tietokantaan, ja n\u00e4in saadaan k\u00e4sitys siit\u00e4 mink\u00e4 verran taivaan t\u00e4hdist\u00e4 on miss\u00e4kin n\u00e4ht\u00e4viss\u00e4.\"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Osallistut maailmanlaajuiseen*\") {\n        $r = $p.Range\n        # Exclude the trailing paragraph-mark character so the paragraph itself\n        # (and its pPr/formatting) is preserved while all runs inside it are\n        # removed.\n        $r.MoveEnd(1, -1)\n        $r.Delete()\n        $r.InsertAfter($newText)\n        break\n    }\n}\n"}
